# BuildingList.xlsx - "Made building models for furnace, archery target, coffin.."
#
# A new tracking header row is inserted at the top of the sheet, pushing all
# existing data down by one row. Three new "progress tracking" columns are
# introduced:
#   D = "Placeholder"       (marked with an "x" for rows that are still
#                             placeholders: Chair, Bed, Table, Coffin)
#   E = "Final"              (header only, no rows marked yet)
#   F = "Needs more coding"  (marked with an "x" for the FarmPlot row)
#
# Furnace (and everything below it) is left unmarked - i.e. it's in good
# shape - matching the commit message about finishing furnace/archery
# target/coffin building models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row down by one to make room for the new header row.
$ws.Rows.Item(1).Insert()

# New header row (row 1).
$ws.Range("D1").Value = "Placeholder"
$ws.Range("E1").Value = "Final"

# Mark the still-placeholder rows in column D.
# (rows 2-7 after the insert correspond to: NONE, Chair, Bed, Table, Coffin, FarmPlot)
$ws.Range("D3").Value = "x"   # Chair
$ws.Range("D4").Value = "x"   # Bed
$ws.Range("D5").Value = "x"   # Table
$ws.Range("D6").Value = "x"   # Coffin

# FarmPlot still needs more coding - flagged in column F instead.
$ws.Range("F7").Value = "x"

# Header for column F (added after the "x" above so shared-string order
# matches: Placeholder, Final, x, Needs more coding).
$ws.Range("F1").Value = "Needs more coding"

# Best-fit the new columns' widths like Excel would after typing the data.
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(6).AutoFit()

# Leave the cursor where the author last left it.
$ws.Range("F8").Select() | Out-Null
